$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.895.08"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.804.69"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.02"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4409"
$ws.Range("E7").Value = "  +4.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3701"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07456"
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8593"
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.69"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.798.77"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.642"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.26"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07069"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.274"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008704"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.80"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "26.911.27"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.167"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.83"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.988"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.28"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.212"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.206"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.35"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08775"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7424"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.164"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.475"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.891"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.096"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01970"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5253"
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.070"
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.822"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1681"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.500"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.21"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9999"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.663"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06337"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9241"
$ws.Range("E51").Value = "  +2.67%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.060"
$ws.Range("E44").Value = "  +6.86%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4987"
$ws.Range("E45").Value = "  +5.98%  "
$ws.Range("E46").Value = "  -1.44%  "
